# Add data for 2021-12-15
# Update the "through" date in the sheet name and the December label,
# and update the December / Total rows with the latest carjacking figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2021-12-07"

# Update the December row label text
$ws.Range("A13").Value = "December (through 12-07)"

# Update 2021 column (H) prior-year-total cell for November row (row 12)
$ws.Range("H12").Value = 200

# Update December row (row 13) values for years 2015-2021 (columns B-H)
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = 22
$ws.Range("D13").Value = 25
$ws.Range("E13").Value = 14
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 35
$ws.Range("H13").Value = 57

# Update Total row (row 14) values for years 2015-2021 (columns B-H)
$ws.Range("B14").Value = 297
$ws.Range("C14").Value = 585
$ws.Range("D14").Value = 846
$ws.Range("E14").Value = 696
$ws.Range("F14").Value = 542
$ws.Range("G14").Value = 1299
$ws.Range("H14").Value = 1699
